$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "Ahmed Mahmoud"
$ws.Range("C3").Value = "Cairo, Egypt"
$ws.Range("D3").Value = "Fri, Sep 29, 2023"
$ws.Range("E3").Value = "3:59 AM"

# Row 4
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Farah Mahmoud"
$ws.Range("C4").Value = "Cairo, Egypt"
$ws.Range("D4").Value = "Fri, Sep 29, 2023"
$ws.Range("E4").Value = "1:02 AM"

# Row 5
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "Farah Mahmoud"
$ws.Range("C5").Value = "Cairo, Egypt"
$ws.Range("D5").Value = "Fri, Sep 29, 2023"
$ws.Range("E5").Value = "4:00 AM"

# Copy the formatting from the existing data row (row 2) onto the new rows
$ws.Range("A2:E2").Copy()
$ws.Range("A3:E5").PasteSpecial(-4122)
